$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Cluster name" column (old column B). This shifts the
# "Number of Nodes (German Version)" and "Number of Nodes (English Version)"
# columns left by one (old C -> B, old D -> C).
$ws.Columns.Item(2).Delete() | Out-Null

# --- Header row ---
$ws.Range("A1").Value = "Cluster"
$ws.Range("B1").Value = "# Nodes German"
$ws.Range("C1").Value = "# Nodes English"

# --- Data rows: merge old cluster number + cluster name into column A ---
$ws.Range("A2").Value = "1- Internet of things & manufacturing"
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 28

$ws.Range("A3").Value = "2- Cognition "
$ws.Range("B3").Value = 23
$ws.Range("C3").Value = 7

$ws.Range("A4").Value = "3- Data, processing & analytics"
$ws.Range("B4").Value = 23
$ws.Range("C4").Value = 19

$ws.Range("A5").Value = "4- Automation  "
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 20

$ws.Range("A6").Value = "5- Cloud computing"
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = 21

$ws.Range("A7").Value = "6- Production, product life circle and flow production"
$ws.Range("B7").Value = 27
$ws.Range("C7").Value = "-"

# --- Column widths for the node-count columns ---
$ws.Columns.Item(2).ColumnWidth = 11.4987
$ws.Columns.Item(3).ColumnWidth = 11.3307

# --- Row 1 no longer has an explicit custom height ---
$ws.Rows.Item(1).AutoFit() | Out-Null

# --- Row heights that change because the wrapped text is different now ---
$ws.Rows.Item(5).RowHeight = 28
$ws.Rows.Item(7).RowHeight = 70

# --- Selection moves to E2 ---
$ws.Range("E2").Select() | Out-Null
